# Weekly update for "Hortaliza, Terminal La Palmera de La Serena - Jengibre":
# a new daily/weekly record is inserted at row 48, pushing the existing
# rows 48-148 down to rows 49-149 (dimension grows from A1:R148 to A1:R149).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 48, shifting rows 48..148 down.
$ws.Rows.Item(48).Insert()

# Fill in the new record (same constant columns as the rest of the sheet).
$ws.Cells.Item(48, 1).Value  = 8
$ws.Cells.Item(48, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(48, 3).Value  = "Coquimbo"
$ws.Cells.Item(48, 4).Value  = 45114
$ws.Cells.Item(48, 5).Value  = 4
$ws.Cells.Item(48, 6).Value  = 100114007
$ws.Cells.Item(48, 7).Value  = "Jengibre"
$ws.Cells.Item(48, 8).Value  = "Sin especificar"
$ws.Cells.Item(48, 9).Value  = "Primera"
$ws.Cells.Item(48, 10).Value = 340
$ws.Cells.Item(48, 11).Value = 17000
$ws.Cells.Item(48, 12).Value = 18000
$ws.Cells.Item(48, 13).Value = 17500
$ws.Cells.Item(48, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(48, 15).Value = "Perú"
$ws.Cells.Item(48, 16).Value = 1346
$ws.Cells.Item(48, 17).Value = 13
$ws.Cells.Item(48, 18).Value = "Hortaliza"
